$d = $word.ActiveDocument

# Locate the paragraph containing "There are formats provided for both TD
# and AC in the report format pdf. Please use it as a reference." -
# the three new paragraphs must be inserted right after it (and before
# the existing blank ListParagraph paragraph that currently follows it).
$anchorText = "There are formats provided for both TD and AC in the report format pdf. Please use it as a reference.`r"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq $anchorText) {
        $target = $p
        break
    }
}

$nextPara = $target.Next()
$r = $nextPara.Range
$r.Collapse(1)  # wdCollapseStart: position right before the existing blank paragraph

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$apos = [char]0x2019

# Build the OOXML for the three new paragraphs, followed by a blank
# ListParagraph paragraph that takes the place of (and keeps empty) the
# paragraph that used to directly follow the anchor paragraph.
$xml = "<w:p $wNs><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"
$xml += "<w:p $wNs><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Another task is </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Business Logic Constraints</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr><w:r><w:t xml:space='preserve'>Since </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>it${apos}s</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'> format is a list of bullet points of constraints, we can all work independently on that at the same time and organize them together for the final report.</w:t></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"

$r.InsertXML($xml)
